$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1179.5
$ws.Range("I18").Value = 1179.5
$ws.Range("K18").Value = 1179.5
$ws.Range("M18").Value = -895.5
$ws.Range("H48").Value = 2082.6667
$ws.Range("J48").Value = 2624
$ws.Range("L48").Value = 7872
$ws.Range("N48").Value = -8456
$ws.Range("H56").Value = 2082.6667
$ws.Range("J56").Value = 2624
$ws.Range("L56").Value = 7872
$ws.Range("N56").Value = -8940
$ws.Range("H86").Value = 22683.572
$ws.Range("I86").Value = 15768.8
$ws.Range("J86").Value = 39970.5
$ws.Range("K86").Value = 15768.8
$ws.Range("L86").Value = 39970.5
$ws.Range("M86").Value = -14645.8
$ws.Range("N86").Value = -42216.5
$ws.Range("H89").Value = 22683.572
$ws.Range("I89").Value = 15768.8
$ws.Range("J89").Value = 39970.5
$ws.Range("K89").Value = 78844
$ws.Range("L89").Value = 199852.5
$ws.Range("M89").Value = -73228
$ws.Range("N89").Value = -211084.5
$ws.Range("H112").Value = 2153
$ws.Range("I112").Value = 716.6667
$ws.Range("J112").Value = 2318.7307
$ws.Range("K112").Value = 2150.0001
$ws.Range("L112").Value = 6956.1921
$ws.Range("M112").Value = -1042.0001
$ws.Range("N112").Value = -9172.1921
$ws.Range("H116").Value = 4730.381
$ws.Range("I116").Value = 4737.4375
$ws.Range("J116").Value = 4707.8
$ws.Range("K116").Value = 4737.4375
$ws.Range("L116").Value = 4707.8
$ws.Range("M116").Value = -1295.4375
$ws.Range("N116").Value = -11591.8
$ws.Range("H126").Value = 59366.668
$ws.Range("J126").Value = 59366.668
$ws.Range("L126").Value = 59366.668
$ws.Range("N126").Value = -69246.66800000001
$ws.Range("H137").Value = 7806.1665
$ws.Range("I137").Value = 1705.1818
$ws.Range("J137").Value = 17393.428
$ws.Range("K137").Value = 5115.5454
$ws.Range("L137").Value = 52180.284
$ws.Range("M137").Value = -2565.5454
$ws.Range("N137").Value = -57280.284
$ws.Range("H140").Value = 93390
$ws.Range("J140").Value = 93390
$ws.Range("L140").Value = 93390
$ws.Range("N140").Value = -103750
$ws.Range("H141").Value = 4320.2856
$ws.Range("J141").Value = 4250
$ws.Range("L141").Value = 12750
$ws.Range("N141").Value = -23110

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6239.9067
$ws.Range("I32").Value = 4584.175
$ws.Range("J32").Value = 28316.334
$ws.Range("K32").Value = 4584.175
$ws.Range("L32").Value = 28316.334
$ws.Range("M32").Value = -4297.175
$ws.Range("N32").Value = -28890.334
$ws.Range("H45").Value = 4433.3335
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623
$ws.Range("H74").Value = 41895.1
$ws.Range("I74").Value = 1731.25
$ws.Range("K74").Value = 1731.25
$ws.Range("M74").Value = -857.25
$ws.Range("H77").Value = 41895.1
$ws.Range("I77").Value = 1731.25
$ws.Range("K77").Value = 8656.25
$ws.Range("M77").Value = -4288.25
$ws.Range("H97").Value = 3488.1785
$ws.Range("I97").Value = 1457.3636
$ws.Range("K97").Value = 1457.3636
$ws.Range("M97").Value = -961.3635999999999
$ws.Range("H110").Value = 3245.8718
$ws.Range("I110").Value = 2149.2188
$ws.Range("K110").Value = 2149.2188
$ws.Range("M110").Value = -104.2188000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14500.181
$ws.Range("I20").Value = 4597.353
$ws.Range("J20").Value = 26970.408
$ws.Range("K20").Value = 4597.353
$ws.Range("L20").Value = 26970.408
$ws.Range("M20").Value = -4350.353
$ws.Range("N20").Value = -27464.408
$ws.Range("H94").Value = 3550.4666
$ws.Range("I94").Value = 3705.2068
$ws.Range("J94").Value = 3270
$ws.Range("K94").Value = 3705.2068
$ws.Range("L94").Value = 3270
$ws.Range("M94").Value = -3254.2068
$ws.Range("N94").Value = -4172

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17437.121
$ws.Range("I31").Value = 7160.7896
$ws.Range("J31").Value = 31383.572
$ws.Range("K31").Value = 7160.7896
$ws.Range("L31").Value = 31383.572
$ws.Range("M31").Value = -6865.7896
$ws.Range("N31").Value = -31973.572
$ws.Range("H34").Value = 17437.121
$ws.Range("I34").Value = 7160.7896
$ws.Range("J34").Value = 31383.572
$ws.Range("K34").Value = 7160.7896
$ws.Range("L34").Value = 31383.572
$ws.Range("M34").Value = -6958.7896
$ws.Range("N34").Value = -31787.572

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1662
$ws.Range("I7").Value = 258.625
$ws.Range("J7").Value = 3533.1667
$ws.Range("K7").Value = 775.875
$ws.Range("L7").Value = 10599.5001
$ws.Range("M7").Value = -663.875
$ws.Range("N7").Value = -10823.5001
$ws.Range("H15").Value = 891.6667
$ws.Range("J15").Value = 1060
$ws.Range("L15").Value = 3180
$ws.Range("N15").Value = -3460
$ws.Range("H132").Value = 2646.2856
$ws.Range("I132").Value = 2341.5
$ws.Range("K132").Value = 21073.5
$ws.Range("M132").Value = -18543.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 31500
$ws.Range("I40").Value = 23000
$ws.Range("K40").Value = 23000
$ws.Range("M40").Value = -22849
$ws.Range("H46").Value = 5100
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 13500
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 13500
$ws.Range("M46").Value = -744
$ws.Range("N46").Value = -13812
$ws.Range("H52").Value = 25203
$ws.Range("J52").Value = 25203
$ws.Range("L52").Value = 25203
$ws.Range("N52").Value = -25721
$ws.Range("H70").Value = 18029.875
$ws.Range("I70").Value = 6931.3335
$ws.Range("K70").Value = 6931.3335
$ws.Range("M70").Value = -6661.3335
$ws.Range("H73").Value = 18029.875
$ws.Range("I73").Value = 6931.3335
$ws.Range("K73").Value = 6931.3335
$ws.Range("M73").Value = -5995.3335
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H126").Value = 11891.8
$ws.Range("I126").Value = 13908.1
$ws.Range("K126").Value = 41724.3
$ws.Range("M126").Value = -39254.3
$ws.Range("H132").Value = 38843.332
$ws.Range("I132").Value = 26609.2
$ws.Range("K132").Value = 79827.60000000001
$ws.Range("M132").Value = -77297.60000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 80692.336
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H88").Value = 39999.5
$ws.Range("J88").Value = 39999.5
$ws.Range("L88").Value = 39999.5
$ws.Range("N88").Value = -40855.5
$ws.Range("H90").Value = 80692.336
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H91").Value = 39999.5
$ws.Range("J91").Value = 39999.5
$ws.Range("L91").Value = 39999.5
$ws.Range("N91").Value = -42963.5
$ws.Range("H136").Value = 17200.207
$ws.Range("I136").Value = 19756.166
$ws.Range("J136").Value = 15396
$ws.Range("K136").Value = 59268.49800000001
$ws.Range("L136").Value = 46188
$ws.Range("M136").Value = -56718.49800000001
$ws.Range("N136").Value = -51288

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11508.56
$ws.Range("I136").Value = 2901.818
$ws.Range("J136").Value = 18271
$ws.Range("K136").Value = 8705.454000000002
$ws.Range("L136").Value = 54813
$ws.Range("M136").Value = -6155.454000000002
$ws.Range("N136").Value = -59913
